# Generate Report for handoff
#
# The "f670f724-...md" file was just handed off again: its status moves
# from "Ready for handoff" to "In Translation" on the Overview sheet and
# on each language sheet, with a fresh "Latest Handoff Datetime" stamped
# for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

# --- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B2").Value = "In Translation"
$ws.Range("D3").Value = "2016-01-26 08:51:45"

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B2").Value = "In Translation"
$ws.Range("D3").Value = "2016-01-26 08:51:56"
